$d = $word.ActiveDocument

$replacements = @(
    @{old="638×5=3190"; new="546×4=2184"},
    @{old="114×6=684"; new="807×3=2421"},
    @{old="687×9=6183"; new="607×8=4856"},
    @{old="534×9=4806"; new="783×2=1566"},
    @{old="629×2=1258"; new="429×5=2145"},
    @{old="797×3=2391"; new="591×4=2364"},
    @{old="119×3=357"; new="361×3=1083"},
    @{old="146×7=1022"; new="458×2=916"},
    @{old="493×2=986"; new="351×9=3159"},
    @{old="301×2=602"; new="220×9=1980"},
    @{old="989×3=2967"; new="567×5=2835"},
    @{old="743×4=2972"; new="662×3=1986"},
    @{old="397×2=794"; new="293×2=586"},
    @{old="926×4=3704"; new="400×8=3200"},
    @{old="891×6=5346"; new="954×9=8586"},
    @{old="919×5=4595"; new="978×3=2934"},
    @{old="515×6=3090"; new="309×3=927"},
    @{old="932×6=5592"; new="921×8=7368"},
    @{old="180×3=540"; new="519×5=2595"},
    @{old="612×2=1224"; new="153×7=1071"},
    @{old="996×4=3984"; new="179×4=716"},
    @{old="806×2=1612"; new="564×6=3384"},
    @{old="560×9=5040"; new="825×9=7425"},
    @{old="102×2=204"; new="489×6=2934"},
    @{old="518×5=2590"; new="611×3=1833"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
